$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new "2022-Q3" row at the top of the
#    data table (row 2), pushing the existing quarters down by one row.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Shift existing data rows (2..8) down to (3..9), bottom-up so we never
# clobber a row before it has been read. Column A is a simple 0-based
# sequence number tied to the row position, so it does not need copying -
# it is rewritten below.
for ($r = 8; $r -ge 2; $r--) {
    $dest = $r + 1
    $wsTotal.Cells.Item($dest, 2).Value = $wsTotal.Cells.Item($r, 2).Value()
    $wsTotal.Cells.Item($dest, 3).Value = $wsTotal.Cells.Item($r, 3).Value()
    $wsTotal.Cells.Item($dest, 4).Value = $wsTotal.Cells.Item($r, 4).Value()
}

# Row 9 is brand new - give its index cell (A9) the same style as the rest
# of the column (copied from A8) before writing its value.
$wsTotal.Range("A8").Copy()
$wsTotal.Range("A9").PasteSpecial(-4122)
$wsTotal.Cells.Item(9, 1).Value = 7

# New top row: 2022-Q3 data.
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 1
$wsTotal.Cells.Item(2, 4).Value = 0.73

# ---------------------------------------------------------------------
# 2) Insert a new "2022-Q3" worksheet right after "总计" and before
#    "2022-Q2", seeded from the "2022-Q2" sheet (same columns/format),
#    then overwrite it with the new quarter's figures.
# ---------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item(2)
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

$wsQ3.Cells.Item(2, 4).Value = "10.53"
$wsQ3.Cells.Item(2, 5).Value = "96.33"
$wsQ3.Cells.Item(2, 6).Value = "6.97"
$wsQ3.Cells.Item(2, 7).Value = "0.7339"
$wsQ3.Cells.Item(2, 8).Value = 7
